# Fix CV hyperlinks so they no longer "auto-download": convert the two
# plain w:hyperlink relationships ("Baltimore Ceasefire" and
# "Opiate crisis") into HYPERLINK field codes (fldChar begin/instrText/
# separate/end), matching the pattern already used elsewhere in the
# document, and relocate the stray _GoBack bookmark to sit after the
# newly-rebuilt "Baltimore Ceasefire" field.

$d = $word.ActiveDocument

# Helper: find a hyperlink in the document by its visible text.
function Find-HyperlinkByText {
    param($Doc, $DisplayText)

    for ($i = 1; $i -le $Doc.Hyperlinks.Count; $i++) {
        $candidate = $Doc.Hyperlinks.Item($i)
        if ($candidate.TextToDisplay -eq $DisplayText) {
            return $candidate
        }
    }
    return $null
}

# Helper: replace the contents of a hyperlink's Range with an equivalent
# "complex" HYPERLINK field (fldChar begin -> instrText -> fldChar
# separate -> display run(s) -> fldChar end), optionally followed by the
# _GoBack bookmark. (Positional params only -- this host's PowerShell
# subset doesn't bind named -Param args reliably.)
function Convert-HyperlinkToField {
    param($Doc, $DisplayText, $InstrText, $WithGoBackBookmark)

    $h = Find-HyperlinkByText $Doc $DisplayText
    $rng = $Doc.Range($h.Range.Start, $h.Range.End)

    $rPr = '<w:rPr><w:rStyle w:val="InternetLink"/><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'

    $body = ''
    $body += '<w:r><w:fldChar w:fldCharType="begin"/></w:r>'
    $body += '<w:r><w:instrText xml:space="preserve">' + $InstrText + '</w:instrText></w:r>'
    $body += '<w:r><w:fldChar w:fldCharType="separate"/></w:r>'
    $body += '<w:r>' + $rPr + '<w:t>' + $DisplayText + '</w:t></w:r>'
    $body += '<w:r>' + $rPr + '<w:fldChar w:fldCharType="end"/></w:r>'
    if ($WithGoBackBookmark) {
        $body += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body><w:p>' + $body + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xmlFrag)
}

# 1) Drop the old stray _GoBack bookmark that sat after "Serious Mental
#    Illness track" -- it is being relocated below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Convert "Opiate crisis" first (it comes after "Baltimore Ceasefire"
#    in the document, so rewriting it first keeps the earlier
#    hyperlink's Range indices valid).
Convert-HyperlinkToField $d 'Opiate crisis' ' HYPERLINK "https://github.com/peterphalen/code-for-publications/tree/master/Phalen-Ray-Watson-Huynh-Greene" \h ' $false

# 3) Convert "Baltimore Ceasefire", re-adding the _GoBack bookmark right
#    after it.
Convert-HyperlinkToField $d 'Baltimore Ceasefire' 'HYPERLINK "https://github.com/peterphalen/code-for-publications/tree/master/Phalen-Bridgeford-Gant-Kivisto-Ray-Fitzgerald" \h ' $true
